$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$ws = $wb.Worksheets.Item("Test")

# Update a few data rows that changed values.
$ws.Range("D3").Value = 1
$ws.Range("F3").Value = "green"
$ws.Range("F7").Value = "blue"
$ws.Range("F10").Value = "red"

# Move the active selection from D10 to D4 (matches the saved view state).
$ws.Range("D4").Select() | Out-Null

# The second worksheet ("Sheet1") is no longer needed - remove it.
$ws2 = $wb.Worksheets.Item("Sheet1")
$ws2.Delete() | Out-Null
